$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 33
$wsMeans.Range("D9").Value = 67
$wsMeans.Range("E9").Value = 68
$wsMeans.Range("F9").Value = 64
$wsMeans.Range("G9").Value = 54

$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.37
$wsMeans.Range("D10").Value = 0.4
$wsMeans.Range("E10").Value = 0.4
$wsMeans.Range("F10").Value = 0.4
$wsMeans.Range("G10").Value = 0.41

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 9.5
$wsSD.Range("D9").Value = 7.1
$wsSD.Range("E9").Value = 5.8
$wsSD.Range("F9").Value = 4.9
$wsSD.Range("G9").Value = 8.3

$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.083
$wsSD.Range("D10").Value = 0
$wsSD.Range("E10").Value = 0.000000000000000045
$wsSD.Range("F10").Value = 0.024
$wsSD.Range("G10").Value = 0.051
